$wb = $excel.ActiveWorkbook

# --- Productdata sheet: scale InventoryCosts (D), BackorderCosts (F)
# and LostSale (I) columns for rows 2-11 by 0.0004 (averaging demand
# figures together with the safety-stock scenario) ---
$ws = $wb.Worksheets.Item("Productdata")

$ws.Cells.Item(2, 4).Value2 = 0.00448
$ws.Cells.Item(2, 6).Value2 = 0.00896
$ws.Cells.Item(2, 9).Value2 = 0.08959999999999999

$ws.Cells.Item(3, 4).Value2 = 0.00496
$ws.Cells.Item(3, 6).Value2 = 0.00992
$ws.Cells.Item(3, 9).Value2 = 0.0992

$ws.Cells.Item(4, 4).Value2 = 0.004920000000000001
$ws.Cells.Item(4, 6).Value2 = 0.009840000000000002
$ws.Cells.Item(4, 9).Value2 = 0.09840000000000002

$ws.Cells.Item(5, 4).Value2 = 0.00444
$ws.Cells.Item(5, 6).Value2 = 0.00888
$ws.Cells.Item(5, 9).Value2 = 0.0888

$ws.Cells.Item(6, 4).Value2 = 0.00048
$ws.Cells.Item(6, 6).Value2 = 0.00096
$ws.Cells.Item(6, 9).Value2 = 0.009600000000000001

$ws.Cells.Item(7, 4).Value2 = 0.00048
$ws.Cells.Item(7, 6).Value2 = 0.00096
$ws.Cells.Item(7, 9).Value2 = 0.009600000000000001

$ws.Cells.Item(8, 4).Value2 = 0.0004400000000000001
$ws.Cells.Item(8, 6).Value2 = 0.0008800000000000001
$ws.Cells.Item(8, 9).Value2 = 0.008800000000000002

$ws.Cells.Item(9, 4).Value2 = 0.00004
$ws.Cells.Item(9, 6).Value2 = 0.00008
$ws.Cells.Item(9, 9).Value2 = 0.0008

$ws.Cells.Item(10, 4).Value2 = 0.00004
$ws.Cells.Item(10, 6).Value2 = 0.00008
$ws.Cells.Item(10, 9).Value2 = 0.0008

$ws.Cells.Item(11, 4).Value2 = 0.00004
$ws.Cells.Item(11, 6).Value2 = 0.00008
$ws.Cells.Item(11, 9).Value2 = 0.0008

# --- ForcastedStandardDeviation sheet: zero out columns B-E for rows 9-11 ---
$ws2 = $wb.Worksheets.Item("ForcastedStandardDeviation")

for ($r = 9; $r -le 11; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $ws2.Cells.Item($r, $c).Value2 = 0
    }
}
